$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. AlbumManager (sheet1): the two test rows for DeleteItem/MoveItem are
#    removed; the old "Save" row (row 4) becomes row 2.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("AlbumManager")
$ws1.Rows("2:3").Delete()

# ---------------------------------------------------------------------------
# 2. Add the new "A-Grade" worksheet after the existing sheets.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "A-Grade"

# Column widths matching the authored layout.
$ws3.Columns(1).ColumnWidth = 35.333333333333336
$ws3.Columns(2).ColumnWidth = 25
$ws3.Columns(3).ColumnWidth = 32.666666666666664
$ws3.Columns(4).ColumnWidth = 32

# Header row (bold "Normal" style, same shared strings as the other sheets).
$ws3.Range("A1").Value = "Method Name"
$ws3.Range("B1").Value = "Purpose"
$ws3.Range("C1").Value = "Test result"
$ws3.Range("D1").Value = "Action"
$ws3.Range("A1:D1").Font.Bold = $true

# Row 2: BugViewModel.Validate
$ws3.Range("A2").Value = "BugViewModel.Validate"
$ws3.Range("B2").Value = "I wanted to test my Validate logic"
$ws3.Range("C2").Value = "The result was to yet again make my private methods internal to be able to tests them. I also realized that both Save and Validate used a class property instead of an in-parameter, which would be more easily tested."
$ws3.Range("B2:C2").WrapText = $true
$ws3.Range("D2").Style = "Dålig"
$ws3.Rows(2).RowHeight = 105

# Row 3: BugViewModel.Save
$ws3.Range("A3").Value = "BugViewModel.Save"
$ws3.Range("B3").Value = "I wanted to test that Validate logic was ok and that the delegate was called"
$ws3.Range("C3").Value = "It failed because when calling OnSave, there were no delegate registered"
$ws3.Range("D3").Value = "Add check that OnSave is not null before invoking the event and also mock the event to prevent it from executing"
$ws3.Range("B3:D3").WrapText = $true
$ws3.Rows(3).RowHeight = 60

# ---------------------------------------------------------------------------
# 3. Restore / update the selection on each sheet (also drives tabSelected /
#    activeTab bookkeeping - the last Select() wins as the active sheet).
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("AlbumService")
$ws2.Range("A1:D1").Select() | Out-Null

$ws1.Range("B3").Select() | Out-Null

$ws3.Range("D3").Select() | Out-Null
